# Se agrego openxlsx a la lista de paquetes
# Add "openxlsx" as a new entry at the end of the package list (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty cell below the existing list in column A and write
# the new package name there (this becomes row 41, right after "corrgram").
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$targetCell = $ws.Cells.Item($newRow, 1)
$targetCell.Value = "openxlsx"
$targetCell.Select() | Out-Null
